$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (columns D, J, K, L, M, O, P) — the rest of the
# columns (A, B, C, E, F, G, H, I, N, Q, R) are identical across rows 2-16
# and are left untouched.
$data = @(
    @{ Row = 2;  D = 44650; J = 130; K = 3000; L = 3500; M = 3308; O = "Región Metropolitana";   P = 551 }
    @{ Row = 3;  D = 44685; J = 150; K = 3000; L = 3500; M = 3267; O = "Región Metropolitana";   P = 544 }
    @{ Row = 4;  D = 44672; J = 140; K = 3000; L = 3500; M = 3286; O = "Región Metropolitana";   P = 548 }
    @{ Row = 5;  D = 44658; J = 180; K = 2500; L = 3000; M = 2778; O = "Región Metropolitana";   P = 463 }
    @{ Row = 6;  D = 44876; J = 80;  K = 6500; L = 7000; M = 6812; O = "Región Metropolitana";   P = 1135 }
    @{ Row = 7;  D = 44671; J = 150; K = 3500; L = 4000; M = 3733; O = "Región Metropolitana";   P = 622 }
    @{ Row = 8;  D = 45021; J = 50;  K = 4500; L = 5000; M = 4700; O = "Región Metropolitana";   P = 783 }
    @{ Row = 9;  D = 44659; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    @{ Row = 10; D = 44643; J = 90;  K = 2800; L = 3000; M = 2911; O = "Región Metropolitana";   P = 485 }
    @{ Row = 11; D = 44987; J = 130; K = 4500; L = 5000; M = 4692; O = "Región Metropolitana";   P = 782 }
    @{ Row = 12; D = 44631; J = 110; K = 3000; L = 3500; M = 3273; O = "Provincia de Chacabuco"; P = 546 }
    @{ Row = 13; D = 44644; J = 140; K = 2500; L = 3000; M = 2786; O = "Provincia de Chacabuco"; P = 464 }
    @{ Row = 14; D = 44630; J = 90;  K = 2500; L = 3000; M = 2722; O = "Región Metropolitana";   P = 454 }
    @{ Row = 15; D = 44637; J = 170; K = 2800; L = 3000; M = 2906; O = "Región Metropolitana";   P = 484 }
    @{ Row = 16; D = 44957; J = 70;  K = 1500; L = 2000; M = 1857; O = "Región Metropolitana";   P = 310 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("J$r").Value = $entry.J
    $ws.Range("K$r").Value = $entry.K
    $ws.Range("L$r").Value = $entry.L
    $ws.Range("M$r").Value = $entry.M
    $ws.Range("O$r").Value = $entry.O
    $ws.Range("P$r").Value = $entry.P
}
